$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handoff"
#
# The localization status report is refreshed: the files that were
# previously "Handed back: in sync with en-US" / "In Translation" are now
# all "Ready for handoff", and the handoff timestamps are updated to the
# new generation time.
# ---------------------------------------------------------------------------

$newStatus = "Ready for handoff"

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overviewDate = "2017-11-06 02:56:39"

$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = $overviewDate

$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Range("G3").Value = $overviewDate

$overview.Range("E4").Value = $newStatus
$overview.Range("F4").Value = $newStatus
$overview.Range("G4").Value = $overviewDate

# Column widths shrink now that the status text is shorter.
$overview.Columns.Item(5).ColumnWidth = 16.333333333333336
$overview.Columns.Item(6).ColumnWidth = 16.333333333333336

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcnDate = "2017-11-06 02:56:33"

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = $zhcnDate

$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("H3").Value = $zhcnDate

$zhcn.Range("C4").Value = $newStatus
$zhcn.Range("H4").Value = $zhcnDate

$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333336

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dedeDate = "2017-11-06 02:56:39"

$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = $dedeDate

$dede.Range("C3").Value = $newStatus
$dede.Range("H3").Value = $dedeDate

$dede.Range("C4").Value = $newStatus
$dede.Range("H4").Value = $dedeDate

$dede.Columns.Item(3).ColumnWidth = 16.333333333333336
